$wb = $excel.ActiveWorkbook

# Update the "想去人数" (want-to-go count) figures on both the
# "展览" sheet and the combined "全部类型" sheet, each of which
# contain the same rows of data.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 2234
    $ws.Range("F6").Value = 780
}
